$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update n_ratings (column I) values per diff
$ws.Range("I2").Value = 3
$ws.Range("I3").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 2
